$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Activate()

# Insert a new (blank) column before column N, shifting the existing
# "Late", "heading" and "Outstanding" columns one place to the right.
$ws.Columns("N").Insert()

# The newly inserted column picks up the width of the column to its
# left (column M).
$ws.Columns("N").ColumnWidth = 9.877604166666666

# Update the selected cell to match the saved view state.
$ws.Range("M16").Select() | Out-Null
